# Add a new beverage record as row 2 of the active sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "53fad91e-c4c8-42f5-a81f-f809f838c37f"
$ws.Range("B2").Value = "dayly"
$ws.Range("C2").Value = "Juices"
$ws.Range("D2").Value = "Coca-Cola"
$ws.Range("E2").Value = 4
$ws.Range("F2").Value = 30

# G2 holds a date-looking string ("2024-09-13"); format it as text first so
# Excel stores the literal string instead of auto-converting it to a date
# serial number.
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "2024-09-13"

$ws.Range("H2").Value = "16:45:10"
